$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ParticipantsTab query (row 2, column B) with the new, more
# elaborate Cypher query text. Assigning a brand-new string automatically
# drops the old shared-string entry that is no longer referenced anywhere
# else in the workbook.
$newParticipantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.platform in ['Illumina NextSeq']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

$ws.Range("B2").Value = $newParticipantsQuery

# The new query text wraps across more lines, so the row grows taller to
# keep showing the whole query (186 -> 279 points).
$ws.Rows.Item(2).RowHeight = 279

# Reflect the view state captured in the saved workbook: scrolled down a
# couple of rows with B5 as the active selection.
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
